# Update "想去人数" (interest count) values in column F for both the
# "展览" and "全部类型" worksheets (which carry the same data).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    6  = 260
    10 = 160
    12 = 4950
    14 = 7224
    18 = 592
    21 = 4237
    22 = 1701
    24 = 89
    25 = 2825
    26 = 577
    28 = 189
    29 = 437
    30 = 408
    31 = 430
    32 = 266
    33 = 73
    35 = 1127
    37 = 1323
    38 = 99
    42 = 15
    43 = 41
    45 = 2613
    46 = 675
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
